$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, "sad",      1.062102600000799),
    @(2, "fear",     1.463099300002796),
    @(3, "neutral",  1.678607100009685),
    @(4, "happy",    2.135661699998309),
    @(5, "sad",      2.712933200004045),
    @(6, "happy",    2.911483200005023),
    @(7, "neutral",  4.82185040001059),
    @(8, "happy",    5.059743100006017),
    @(9, "neutral",  5.875170300001628),
    @(10, "angry",   6.073256500007119),
    @(11, "happy",   6.470396800010349),
    @(12, "neutral", 7.110622200008947),
    @(13, "angry",   7.498422500008019),
    @(14, "fear",    8.126322500000242),
    @(15, "surprise",13.64557620001142),
    @(16, "fear",    18.32348210000782),
    @(17, "angry",   18.54596410000522),
    @(18, "fear",    18.74780310000642),
    @(19, "neutral", 19.64774410000246),
    @(20, "sad",     20.6621487000084)
)

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "neutral"
$ws.Range("C2").Value = 0

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 3 + $i
    $entry = $data[$i]
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
}

# Rows 10-22 are brand new rows; column A needs the same style (s="1")
# that rows 2-9 already carry. Copy the format from an existing styled
# cell (A9) and paste only the formatting onto the new A column cells.
$ws.Range("A9").Copy()
$ws.Range("A10:A22").PasteSpecial(-4122)
